$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column A (rows 1-32)
$colA = @(
    -0.17413792331684164,
    -0.15170596666803338,
    -0.10069945673606817,
    -0.09235413559450123,
    -0.089053061229165564,
    -0.04314296658595751,
    -0.032702133097198605,
    0.019495993644183418,
    0.021528638942418077,
    0.023577155866783883,
    -0.024387042357473909,
    -0.020861174461332066,
    -0.017168543588584662,
    -0.0090810182174791976,
    -0.0080525987830668555,
    -0.0060342641636763616,
    -0.0040032669969436441,
    -0.016102792028064528,
    -0.012091153287681067,
    -0.0080162599013053892,
    -0.0040055438207966532,
    -0.045714888667395925,
    -0.040500896170500233,
    -0.020099319876584865,
    -0.097253873361879073,
    -0.094629571426045445,
    -0.091969132007687193,
    -0.089015396838727945,
    -0.081362517668230439,
    -0.021173479813068496,
    -0.014023000494658433,
    -0.0040009213546099431
)

# New values for column B (rows 1-32)
$colB = @(
    0.17380875865077172,
    0.15041434536290677,
    0.10035413552039962,
    0.092053061186891938,
    0.088044057925706909,
    0.042702132992475939,
    0.032599564583661866,
    -0.019528638988909552,
    -0.021577155909847434,
    -0.023582164423762819,
    0.024361174405738595,
    0.020668543528838512,
    0.017081018123939806,
    0.0090525987399106,
    0.0080342641125952241,
    0.0060032669449858744,
    0.0039999999331978842,
    0.016091153256660107,
    0.012016259867849488,
    0.0080055437869717139,
    0.0039999999659006136,
    0.045500896120016066,
    0.040099319712648018,
    0.019999999833713034,
    0.097129571374853896,
    0.094469131953264451,
    0.091015396772316404,
    0.0883625175550522,
    0.081173479308259466,
    0.021023000372334177,
    0.014000921209914807,
    0.0039999998988697882
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# Column B width narrows from 15.42578125 to 14.7109375 (character units ~13.88)
$ws.Columns.Item(2).ColumnWidth = 13.83
